$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: YDS  (B2/C2 = OFF R / DEF R ; B3/C3 = OFF P / DEF P)
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value2 = @"
7 0 4 7 0 3 10 3 6 3 3 1 9 7 1 0 -2 0 25 3 0 1 12 7 0 2 7 -2 6 17 16 3 1 9 -2 -2 2 6 0 2 1 3 0 4 13 4 3 10 2 2 -1 3 9 2 -2 -1 2 3 9 8 0 2 2 5 15 5 -1 2 0 0 3 3 9 4 1 -1 6 1 19 2 4 2 -1 0 2 2 2 0 -3 3 -3 43 4 3 4 4 2 11 9 2 3 20 2 9 3 6 4 0 5 1 9 2 1 4 1 -2 3 -2 8 0 3 -3 1 9 14 -1 0 -1 6 3 2 14 2 6 20 9 4 19 6 1 1 6 5 10 8 -1 10 6 4 2 4 11 0 3 2 3 3 3 6 2 -1 1 5 0 -1 -2 2 55 2 7 7 20 3 10 3 -5 3 2 15 6 1 3 6 0 2 4 4 9 4 1 6 5 16 8 10 4 2 0 3 8 1 6 -4 3 6 7 2 0 -3 4 11 9 2 0 1 20 -2 14 9 1 2 3 -2 5 25 7 5 18 3 2 3 -1 9 4 2 20 1 4 9 20 -2 8 0 3 -2 2 13 3 4 -3 1 1 0 5 1 -1 6 3 0 2 11 -3 0 4 2 -5 2 3 7 0 3 11 0 -1 5 23 1 0 0 3 1 0 2 8 0 4 8 8 4 3 65 1 6 3 3 13 2 0 3 5 2 7 1 3 -1 4 5 1 8 6 1 9 3 3 2 5 0 4 6 5 9 1 4 4 9 -1 0 15 0 2 11 5 4 -1 6 1 -4 2 10 4 13 1 2 3 0 1 14 2 10 2 4 4 0 3 3 14 5 6 7 9 8 8 0 3 2 8 6 8 0 6 3 1 17 4 5 2 1 1 12 3 -1 1 10 5 3 6 9 3 5 6 1 2 6 7 5 3 5 5 0 20 1 -1 6 2 1 -2 5 0 4 0 3 0 1 4 5 2 1 -1 2 10 2 0 7 15 3 7 -1 2 5 2 9 2 2 -2 8 7 3 0 5 4 6 6 2 2 12 3 8 0 2 4 70 3 -1 2 0 5 9 3 2 7 0 1 -6 0 3 3 16 2 2 5 0 5 8 -1 3 14 2 0 3 6 5 7 3 8 4 11 0 1 14 2 2 0 13 4 13 1 10 0 -1 3 2 -3 10 4 7 2 -1 1 -2 1 0 -2 9 2 2 6 -1 7 8 1 4 14 3 2 31 3 3 3 3 -3 6 10 -1 12 1 5 2 -4 8 0 49 2 7 1 1 2 3 10 -1 11 4 5 10 1 9 3 11 3 4 0 3 0 5 8 4 0 30 4 2 4 4 1 1 0 6 2 1 6 0 2 7 3 3 5 2 0 7 2 8 4 2 6 2 2 7 1 1 7 2 -1 5 3 0 11 14 1 7 1 -2 -2 11 20 -3 7 13 1 -2 3 17 0 12 2 2 11 2 2 2 30 2 1 1 7 14 6 7 6 4 -4 2 0 1 10 3 2 0 2 1 -1 1 3 9 11 2 15 7 1 1 4 20 -4 8 9 1 2 12 2 10 2 4 2 4 2 -2 11 12 8 16 -1 9 1 2 2 6 6 2 2 0 3 0 7 4 3 8 4 4 3 9 -1 1 14 2 9 6 7 2 14 3 0 1 11 3 5 5 -1 5 4 8 9 7 2 6 2 16 1 11 5 3 1 13 9 1 5 6 3 3 9 4 1 3 2 2 1 3 3 -1 6 16 5 5 5 9 14 -1 1 1 6 8 3 4 -1 10 5 1
"@

$wsYDS.Range("C2").Value2 = @"
11 4 1 3 10 6 3 2 0 11 2 3 3 2 5 3 1 0 2 1 5 2 8 4 5 2 0 4 0 0 5 13 4 5 0 8 11 2 2 3 12 2 2 3 -1 2 -2 3 4 -2 5 59 4 0 0 4 -3 2 5 6 2 9 4 4 4 -1 2 1 0 2 13 4 5 4 -1 3 8 3 46 -2 5 0 6 2 3 4 2 2 5 1 4 2 4 4 6 1 8 13 6 2 12 7 1 6 0 1 -1 2 2 0 4 0 7 2 2 38 1 6 1 12 3 13 2 4 11 16 1 2 27 0 1 16 3 1 1 0 0 2 3 1 1 2 2 22 4 1 1 5 0 2 22 3 4 9 6 4 14 5 6 0 6 8 26 1 3 2 3 4 5 -4 7 3 11 1 18 4 3 0 4 13 4 2 1 3 8 4 -3 4 5 8 -4 -2 3 2 6 12 3 3 3 2 1 4 3 0 4 0 1 3 2 3 9 11 4 3 1 1 9 3 5 0 24 2 3 13 6 3 -5 11 2 1 9 9 2 5 1 1 4 5 7 3 1 8 2 9 2 23 4 4 7 2 1 4 3 -3 3 2 14 2 0 4 3 1 5 4 6 3 3 1 4 7 15 4 14 1 5 5 2 6 1 5 7 1 -1 3 2 5 0 1 3 11 8 1 2 17 36 3 0 1 8 3 11 6 1 8 -2 4 4 1 4 2 5 6 30 20 3 0 16 3 1 6 3 10 6 5 12 4 1 -3 1 4 11 14 -2 2 6 10 3 9 5 3 7 12 5 6 4 7 10 2 10 4 8 5 2 6 5 24 4 1 0 6 4 3 9 26 6 1 1 10 2 1 3 51 3 2 4 5 3 4 7 5 0 -1 3 -2 0 5 2 11 23 1 4 1 3 0 2 1 5 6 2 2 18 2 3 16 2 4 3 28 5 3 1 5 -2 4 2 -1 1 5 2 3 7 3 5 2 4 2 6 1 4 4 4 2 0 6 2 2 1 6 4 -1 12 7 11 5 4 1 4 2 2 2 13 8 5 4 2 -7 3 -1 7 5 2 1 1 3 4 5 7 11 2 1 -1 6 -2 3 4 3 5 4 0 2 3 5 3 5 5 4 5 6 3 3 3 5 7 6 5 2 4 5 20 -3 3 5 3 0 0 1 3 13 -2 5 18 3 8 -1 -1 4 9 2 1 4 4 1 4 6 2 6 4 3 9 3 6 18 6 2 0 1 3 1 5 3 1 3 0 20 10 4 12 1 9 2 1 2 5 10 11 2 3 3 2 0 9 1 6 1 9 5 -2 20 6 3 7 3 5 8 0 2 10 4 3 2 2 5 3 6 5 0 4 4 7 0 4 5 13 12 12 4 4 4 5 8 1 7 11 7 2 3 -2 7 2 5 10 9 10 17 -2 0 1 4 3 -8 7 4 2 5 11 31 0 3 4 3 0 5 -1 5 21 -2 15 5 25 1 7 5 3 1 2 6 1 23 1 1 7 5 6 3 3 5 1 5 -1 2 13 5 2 1 13 10 2 3 5 10 4 2 12 4 2 0 5 4 2 4 4 4 4 6 7 -4 1 5 1 13 35 5 9 3 4 4 2 3 5 4 -3 7 3 3 6 -3 8
"@

$wsYDS.Range("B3").Value2 = @"
17 3 12 8 -6 9 6 2 16 31 11 8 25 7 -2 21 -2 15 12 5 18 20 20 19 1 15 45 17 13 14 8 20 -2 1 1 18 17 16 7 22 5 9 6 2 15 18 7 -2 4 15 28 5 9 5 4 4 5 2 22 7 6 3 1 21 6 -2 7 7 48 40 11 9 3 8 6 13 9 7 7 10 7 31 41 7 4 18 27 17 16 35 16 8 3 27 5 2 10 7 19 5 17 13 21 12 13 4 6 6 13 11 24 5 11 9 8 4 5 3 4 11 4 12 12 9 6 18 15 2 43 7 10 9 40 6 2 7 14 3 5 6 1 32 9 18 14 9 7 9 12 9 3 7 2 7 24 3 20 36 4 11 11 9 3 9 41 4 1 9 14 17 0 11 26 27 1 9 27 14 6 6 7 11 13 8 4 9 14 16 15 2 9 15 4 41 9 22 1 15 9 6 9 16 7 0 29 13 37 7 5 8 11 0 10 5 5 8 21 11 10 8 5 3 8 4 13 7 6 2 3 4 2 31 2 11 37 14 32 29 2 20 49 1 10 16 11 -2 8 5 5 9 6 6 19 6 -6 9 -2 6 4 5 10 7 6 5 7 5 17 2 9 10 10 16 10 11 18 7 12 15 10 3 5 15 19 4 7 17 14 17 29 1 16 8 4 1 9 9 9 1 2 11 12 4 20 7 13 10 1 92 25 15 13 7 6 1 9 5 25 -4 5 11 6 2 6 17 16 2 10 6 8 20 7 22 13 3 6 13 14 8 6 6 7 4 10 16 33 15 12 9 55 14 5 10 11 9 0 4 11 14 6 6 10 8 5 28 6 17 6 9 27 6 21 -3 9 7 24 14 4 31 11 19 3 3 12 14 3 11 11 8 3 5 6 13 -3 9 9 32 9 5 10 9 6 14 12 12 3 12 19 10 9 15 9 39 23 6 9 5 23 6 7 14 23 5 1 3 6 11 7 5 7 10 14 3 14 11 6 10 5 7 26 13 12 12 12 26 6 18 16 6 4 4 1 6 6 8 31 1 0 6 6 9 12 13 4 11 9 8 8 5 2 9 10 3 7 9 10 1 31 5 31 12 19 19 15 7 3 4 -4 9 2 11 14 9 18 7 -3 25 5 8 11 44 19 7 6 9 17 6 3 11 20 8 -3 16 40 -2 7 5 64 4 13 15 26 10 7 4 10 6 15 7 5 -9 0 13 9 14 3 3 13 19 5 7 12 7 5 2 12 1 12 14 42 1 36 10 10 6 3 2 34 5 9 5 8 11 9 13 19 12 6 19 10 7 13 32 8 9 8 6 7 7 4 8 7 9 10 16 5 5 4
"@

$wsYDS.Range("C3").Value2 = @"
22 16 20 4 5 1 23 10 13 9 12 4 7 5 7 6 -1 8 7 4 13 7 5 5 8 9 6 11 6 8 15 8 5 18 7 10 18 84 9 6 9 -4 2 4 11 10 2 6 0 21 8 -4 5 7 28 12 3 10 5 9 7 3 9 10 10 47 6 6 1 33 9 1 27 6 17 11 6 7 10 8 11 7 30 9 26 10 6 3 9 5 18 11 5 3 21 27 4 10 3 11 9 8 35 7 5 11 4 10 7 13 9 7 7 12 10 16 14 12 10 16 5 19 -2 22 6 -3 16 6 11 8 19 14 38 4 5 15 17 7 6 17 7 22 10 11 2 3 6 4 8 6 19 4 7 4 7 7 34 22 -7 6 9 11 1 36 2 3 24 5 3 16 5 5 15 13 11 5 5 -4 7 8 6 3 51 7 5 11 42 9 6 12 11 21 6 21 10 24 8 -4 3 25 9 2 3 8 8 1 16 8 24 0 5 21 13 5 6 3 13 8 4 9 3 6 10 8 13 6 7 0 25 11 15 1 8 6 9 21 12 2 2 3 -2 3 8 24 4 20 16 19 14 18 20 1 19 6 4 10 -5 15 28 5 20 12 15 30 14 5 20 8 5 17 9 9 0 5 13 4 8 3 3 8 18 4 0 10 21 12 0 9 19 17 1 2 14 8 23 16 0 4 20 32 8 7 2 15 8 5 3 9 12 4 17 15 11 11 27 8 22 22 10 7 6 7 55 6 3 20 4 16 27 14 3 3 25 10 10 9 5 2 29 14 12 21 6 7 26 -9 11 23 9 23 9 57 5 9 8 28 18 12 5 15 15 26 5 14 7 25 6 15 2 14 37 10 21 8 42 -6 4 17 1 6 8 37 2 15 20 8 17 14 6 16 16 6 9 15 8 6 4 24 7 13 25 -2 -8 6 5 19 0 6 6 13 11 4 8 7 6 5 6 -5 14 5 3 22 7 4 10 23 4 13 13 20 6 12 24 12 9 49 32 7 13 14 4 8 13 7 8 17 4 7 14 32 4 23 50 4 5 13 6 59 5 14 14 3 18 1 22 16 9 48 2 3 25 11 16 6 8 8 10 33 31 29 3 40 8 51 7 34 4 5 14 8 4 6 9 15 6 12 -2 6 17 16 18 1 4 8 7 8 27 17 2 13 -2 20 4 12 1 17 14 7 5 14 4 23 20 19 5 22 14 3 1 5 0 32 17 21 4 18 5 4 24 9 32 5 10 7 4 20 11 4 16 18 3 4 6 24 36 8 6 4 5 9 11 12 14 2 7 24 -2 10 5 9 13 6 18 14 12 9 11 13 17 8 8 4 2 7 13 18 34 9 5 -2 16 9 9 12 13 9 29 14 -7 12 6 14 14 38 8 4 1 13 7 7 2 19 7 5 20 11 -2 9 9 11 13 4 6 14 8 14 7 2
"@

# ---------------------------------------------------------------------------
# Sheet: OFF
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value2 = 350
$wsOFF.Range("D2").Value2 = 16
$wsOFF.Range("E2").Value2 = 27
$wsOFF.Range("F2").Value2 = 135
$wsOFF.Range("G2").Value2 = 97
$wsOFF.Range("J2").Value2 = 42
$wsOFF.Range("L2").Value2 = 480
$wsOFF.Range("M2").Value2 = 304
$wsOFF.Range("O2").Value2 = 37
$wsOFF.Range("P2").Value2 = 18
$wsOFF.Range("Q2").Value2 = 914

$wsOFF.Range("B3").Value2 = 18
$wsOFF.Range("C3").Value2 = 298
$wsOFF.Range("F3").Value2 = 164
$wsOFF.Range("G3").Value2 = 48
$wsOFF.Range("H3").Value2 = 67
$wsOFF.Range("I3").Value2 = 121
$wsOFF.Range("J3").Value2 = 76
$wsOFF.Range("N3").Value2 = 35

# ---------------------------------------------------------------------------
# Sheet: DEF
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value2 = 336
$wsDEF.Range("E2").Value2 = 20
$wsDEF.Range("F2").Value2 = 119
$wsDEF.Range("G2").Value2 = 87
$wsDEF.Range("I2").Value2 = 14
$wsDEF.Range("J2").Value2 = 54
$wsDEF.Range("L2").Value2 = 530
$wsDEF.Range("M2").Value2 = 339
$wsDEF.Range("O2").Value2 = 37
$wsDEF.Range("P2").Value2 = 15
$wsDEF.Range("Q2").Value2 = 959

$wsDEF.Range("C3").Value2 = 308
$wsDEF.Range("F3").Value2 = 171
$wsDEF.Range("G3").Value2 = 63
$wsDEF.Range("H3").Value2 = 56
$wsDEF.Range("I3").Value2 = 90
$wsDEF.Range("J3").Value2 = 120
$wsDEF.Range("N3").Value2 = 45

# ---------------------------------------------------------------------------
# Sheet: ST
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value2 = 135
$wsST.Range("D2").Value2 = 113
$wsST.Range("F2").Value2 = 253
$wsST.Range("G2").Value2 = 246
$wsST.Range("N2").Value2 = 54
$wsST.Range("O2").Value2 = 29

$wsST.Range("B3").Value2 = 105

$wsST.Range("D3").Value2 = @"
51 56 53 52 32 46 35 54 33 48 41 49 47 47 58 41 49 31 46 57 60 44 40 51 46 61 69 58 26 64 34 48 53 44 55 44 48 52 56 53 55 53 47 42 43 39 39 40 36 41 40 38 58 45 52 52 43 56 44 38 42 53 47 19 49 37 41 41 45 48 46 43 53 26 50 47 57 49 39 56 44 28 68 53 52 52 39 33 47 33 46 38 46 59 47 50 44 48 43 41 56 60 44 52 42 43 40 56 56 48 47 41 42
"@

$wsST.Range("D4").Value2 = @"
11 9 -1 0 0 18 0 0 0 0 5 13 3 3 0 0 11 0 0 11 10 0 0 -1 30 0 0 9 0 13 0 0 0 29 0 16 0 6 0 0 19 0 0 0 0 0 0 5 0 0 0 0 13 0 0 7 7 0 0 0 0 21 0 0 0 0 2 0 2 0 0 0 0 0 12 0 0 10 0 42 0 0 0 0 0 0 0 0 0 5 0 0 3 10 0 0 12 0 0 0 0 10 0 9 0 0 0 10 3 0 0 0 -4
"@

$wsST.Range("D5").Value2 = "0 0 2 11"

# ---------------------------------------------------------------------------
# Sheet: TURNS
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("C2").Value2 = 15
$wsTURNS.Range("E2").Value2 = 11
$wsTURNS.Range("D3").Value2 = 13

# ---------------------------------------------------------------------------
# Sheet: PEN
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("B3").Value2 = 23
